$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates (F1,G1,H1 renamed; I1,J1,K1 added)
$ws.Range("F1").Value = "Fit time sktree"
$ws.Range("G1").Value = "Prediction time sktree"
$ws.Range("H1").Value = "Score sktree"
$ws.Range("I1").Value = "Fit time MeanSDTD6"
$ws.Range("J1").Value = "Prediction time MeanSDTD6"
$ws.Range("K1").Value = "Score MeanSDTD6"

# Apply header style (same as A1:H1) to the new header cells
$ws.Range("F1:H1").Copy()
$ws.Range("I1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 (Iris)
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 0.002012252807617188
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0.01124763488769531
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0.8666666666666667

# Row 3 (Wine)
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.9444444444444444
$ws.Range("I3").Value = 0.006361961364746094
$ws.Range("J3").Value = 0.01128768920898438
$ws.Range("K3").Value = 0.9444444444444444

# Row 4 (Breast Cancer)
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 0.01556801795959473
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.9473684210526315
$ws.Range("I4").Value = 0.02974939346313477
$ws.Range("J4").Value = 0.01807379722595215
$ws.Range("K4").Value = 0.9473684210526315

# Row 5 (Digits)
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 0.0401148796081543
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.8527777777777777
$ws.Range("I5").Value = 0.09845972061157227
$ws.Range("J5").Value = 0.08024406433105469
$ws.Range("K5").Value = 0.7555555555555555

# Row 6: Adult -> BankNote Authentication (new data)
$ws.Range("A6").Value = "BankNote Authentication"
$ws.Range("B6").Value = "classification"
$ws.Range("C6").Value = 1372
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.9818181818181818
$ws.Range("I6").Value = 0.04781126976013184
$ws.Range("J6").Value = 0.03228092193603516
$ws.Range("K6").Value = 0.9636363636363636

# Row 7: Gas Drift (new row)
$ws.Range("A7").Value = "Gas Drift"
$ws.Range("B7").Value = "classification"
$ws.Range("C7").Value = 13910
$ws.Range("D7").Value = 128
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 3.521384954452515
$ws.Range("G7").Value = 0.001997709274291992
$ws.Range("H7").Value = 0.9410496046010065
$ws.Range("I7").Value = 0.5315756797790527
$ws.Range("J7").Value = 0.615117073059082
$ws.Range("K7").Value = 0.8461538461538461

# Row 8: Shuttle (new row)
$ws.Range("A8").Value = "Shuttle"
$ws.Range("B8").Value = "classification"
$ws.Range("C8").Value = 58000
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 0.1600716114044189
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.9994827586206897
$ws.Range("I8").Value = 0.1409909725189209
$ws.Range("J8").Value = 2.074449777603149
$ws.Range("K8").Value = 0.9876724137931034
